$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data shifts down to rows 2-127
$ws.Rows.Item(1).Insert()

# Header cells: use shared-string text, plain (non-date) style
$ws.Range("A1:B1").Style = "Normal"
$ws.Range("A1").Value = "ac"
$ws.Range("B1").Value = "date"

# Add the new 127th data row (A=42 group, continuing the pattern)
$ws.Range("A127").Value = 42
$ws.Range("B127").Style = "Normal"
$ws.Range("B127").NumberFormat = "yyyy\-mm\-dd"

# Update column B (date) values for rows 2-127 to the refreshed draft data
$bValues = @(
  45292,
  45962,
  46631,
  45292,
  46054,
  46722,
  45292,
  46054,
  46784,
  45292,
  46082,
  46813,
  45292,
  45962,
  46661,
  45292,
  45962,
  46722,
  45292,
  46054,
  46753,
  45292,
  45931,
  46631,
  45292,
  46082,
  46844,
  45292,
  46023,
  46722,
  45292,
  46113,
  46784,
  45292,
  46023,
  46661,
  45292,
  46054,
  46905,
  45292,
  45962,
  46631,
  45292,
  45992,
  46722,
  45292,
  46082,
  46753,
  45292,
  46023,
  46692,
  45292,
  45901,
  46631,
  45292,
  45992,
  46631,
  45292,
  45901,
  46692,
  45292,
  45992,
  46784,
  45292,
  46082,
  46784,
  45292,
  45962,
  46631,
  45292,
  45992,
  46722,
  45292,
  46082,
  46813,
  45292,
  46023,
  46692,
  45292,
  46054,
  46813,
  45292,
  45992,
  46722,
  45292,
  46054,
  46874,
  45292,
  46023,
  46784,
  45292,
  46054,
  46784,
  45292,
  46054,
  46844,
  45292,
  46023,
  46784,
  45292,
  45992,
  46692,
  45292,
  45931,
  46753,
  45292,
  46023,
  46661,
  45292,
  45992,
  46753,
  45292,
  45962,
  46722,
  45292,
  46023,
  46753,
  45292,
  46143,
  46966,
  45292,
  45992,
  46722,
  45292,
  45931,
  46661
)

for ($i = 0; $i -lt $bValues.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Restore the selection Excel would show after this edit
$ws.Range("A2:B127").Select()
